$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.564.36"
$ws.Range("E2").Value = "  +6.58%  "
$ws.Range("D3").Value = "2.379.25"
$ws.Range("E3").Value = "  +4.35%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "111.54"
$ws.Range("E5").Value = "  +7.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "317.00"
$ws.Range("E6").Value = "  +2.02%  "
$ws.Range("E7").Value = "  +2.56%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("E9").Value = "  +5.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.91"
$ws.Range("E10").Value = "  +8.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0930"
$ws.Range("E11").Value = "  +3.56%  "
$ws.Range("E12").Value = "  +5.18%  "
$ws.Range("E13").Value = "  +4.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.109"
$ws.Range("E14").Value = "  +0.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.73"
$ws.Range("E15").Value = "  +4.95%  "
$ws.Range("D16").Value = "2.739.54"
$ws.Range("E16").Value = "  +4.23%  "
$ws.Range("D17").Value = "2.382.49"
$ws.Range("E17").Value = "  +4.18%  "
$ws.Range("D18").Value = "45.305.71"
$ws.Range("E19").Value = "  +5.39%  "
$ws.Range("E20").Value = "  +3.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.02"
$ws.Range("E21").Value = "  -2.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.04"
$ws.Range("E22").Value = "  +3.03%  "
$ws.Range("E23").Value = "  +3.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.75"
$ws.Range("E24").Value = "  +2.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.33"
$ws.Range("E25").Value = "  +7.36%  "
$ws.Range("E26").Value = "  -0.55%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.30"
$ws.Range("E27").Value = "  +6.22%  "
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.56"
$ws.Range("E28").Value = "  +9.03%  "
$ws.Range("E29").Value = "  +2.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.93"
$ws.Range("E30").Value = "  +3.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "38.74"
$ws.Range("E31").Value = "  +8.76%  "
$ws.Range("E32").Value = "  +10.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "169.84"
$ws.Range("E33").Value = "  +3.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.03"
$ws.Range("E34").Value = "  +18.63%  "
$ws.Range("E35").Value = "  +2.88%  "
$ws.Range("E36").Value = "  +8.92%  "
$ws.Range("E37").Value = "  +4.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.05"
$ws.Range("E38").Value = "  +12.87%  "
$ws.Range("E39").Value = "  +5.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.91"
$ws.Range("E40").Value = "  +5.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.73"
$ws.Range("E41").Value = "  +11.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.24"
$ws.Range("E42").Value = "  +8.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.81"
$ws.Range("E43").Value = "  +16.50%  "
$ws.Range("E44").Value = "  +6.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "71.52"
$ws.Range("E45").Value = "  +4.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "117.92"
$ws.Range("E47").Value = "  +7.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.79"
$ws.Range("E48").Value = "  +11.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.64"
$ws.Range("E49").Value = "  +20.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "79.25"
$ws.Range("E50").Value = "  +3.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.20"
$ws.Range("E51").Value = "  +7.00%  "
